$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the Q-column helper cells (labels + values) BEFORE re-pointing
#    the MAX / MIN defined names, so the dependent formulas pick up the
#    correct final values once the names are (re)created.
# ---------------------------------------------------------------------------

# Q1 label -> "Max period"
$ws.Range("Q1").Value = "Max period"

# Q2: plain literal value (formula removed)
$ws.Range("Q2").ClearContents()
$ws.Range("Q2").Value = 159

# Q4 label -> "Max value to use"
$ws.Range("Q4").Value = "Max value to use"

# Q5: now holds a formula derived from Q2 (this is the new MAX)
$ws.Range("Q5").Formula = "=ROUNDDOWN(Q2*0.95,0)"

# Q7 label -> "Min value to use"
$ws.Range("Q7").Value = "Min value to use"

# Q8: new cell, literal value (this is the new MIN)
$ws.Range("Q8").Value = 0

# ---------------------------------------------------------------------------
# 2. Re-point the MAX / MIN defined names. Deleting + re-adding (instead of
#    just touching RefersTo) forces the calc engine to rebuild the
#    dependency graph so downstream formulas recalculate correctly.
# ---------------------------------------------------------------------------
$wb.Names.Item("MAX").Delete()
$wb.Names.Item("MIN").Delete()
$wb.Names.Add("MAX", "=Sheet1!`$Q`$5")
$wb.Names.Add("MIN", "=Sheet1!`$Q`$8")

# ---------------------------------------------------------------------------
# 3. Column widths: column I (9) becomes its own width-6 band, splitting the
#    old single 2-10 band into 2-8 (5.5703125), 9 (6), 10 (5.5703125).
#    (COM ColumnWidth is in "characters"; the stored OOXML width is
#    ColumnWidth + 5/6, so subtract that offset to land exactly on 6.)
# ---------------------------------------------------------------------------
$ws.Columns(9).ColumnWidth = 5.166666666666667

# ---------------------------------------------------------------------------
# 4. Angle table (rows 2-7): rebuild as a running formula chain.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Formula = "=A2+1.5"
$ws.Range("C2:J2").Formula = "=B2+1.5"

$ws.Range("A3").Formula = "=A2+13.5+1.5"
$ws.Range("B3:J3").Formula = "=A3+1.5"

$ws.Range("A4").Formula = "=A3+13.5+1.5"
$ws.Range("B4:J4").Formula = "=A4+1.5"

$ws.Range("A5").Formula = "=A4+13.5+1.5"
$ws.Range("B5:J5").Formula = "=A5+1.5"

$ws.Range("A6").Formula = "=A5+13.5+1.5"
$ws.Range("B6:J6").Formula = "=A6+1.5"

$ws.Range("A7").Formula = "=A6+13.5+1.5"
$ws.Range("B7:J7").Formula = "=A7+1.5"

# ---------------------------------------------------------------------------
# 5. Force the sine-table formulas (rows 10-15) to recalculate now that
#    MAX / MIN and the angle table have changed.
# ---------------------------------------------------------------------------
$ws.Range("A10:J15").Formula = $ws.Range("A10:J15").Formula
$excel.CalculateFull()

# ---------------------------------------------------------------------------
# 6. Restore the sheet selection.
# ---------------------------------------------------------------------------
$ws.Range("J28").Select()

Write-Host "done"
